$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @()
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.05789840000000268)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.003664900000003968)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.003326200000000057)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.005946800000003805)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.006491099999998085)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.00386599999999504)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.001894900000003474)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.0007946999999859372)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.004845499999987624)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.003720899999990479)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.001389699999990057)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.002035599999999249)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.002667100000010691)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.01978369999999074)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.02319819999999595)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.001183900000000904)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 500, 50, 0.480582499999997)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 500, 50, 0.008589999999998099)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 500, 50, 0.04207550000000992)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, -0.527860600000011)
$data += ,@("best_of_all_selection", "real_values_crossover_heurestic", "real_values_uniform_mutation", 50, 50, 0.0002268000000071879)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
